$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge the two runs "Sat Sep 15" and " 12:06:42 PDT 2017" into a
#    single run "Sat Sep 15 12:06:42 PDT 2017".
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Sat Sep 15 12:06:42 PDT 2017", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sat Sep 15 12:06:42 PDT 2017", 2)

# ---------------------------------------------------------------------
# 2. Append the new "Sun Sep 16" purchase-details block right after the
#    "Amount balance ... - 138337.0" paragraph, before the trailing
#    blank paragraphs that were already at the end of the document.
# ---------------------------------------------------------------------

# Find the (still bold) "Amount balance" paragraph that ends in 138337.0
# and the first (non-bold) blank paragraph that immediately follows it.
$boldAnchor = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*138337.0*") {
        $boldAnchor = $i
        break
    }
}

# Two additional bold, empty paragraphs right after the bold anchor.
$anchorRange = $d.Paragraphs.Item($boldAnchor).Range
$anchorRange.InsertParagraphAfter()
$anchorRange2 = $d.Paragraphs.Item($boldAnchor + 1).Range
$anchorRange2.InsertParagraphAfter()

# Index of the first untouched (plain) trailing paragraph - it has
# shifted two slots further down because of the two inserts above.
$global:plainAnchor = $boldAnchor + 1 + 2

function Insert-PlainParagraph([string]$text) {
    $target = $d.Paragraphs.Item($global:plainAnchor)
    $target.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($global:plainAnchor)
    $newPara.Range.Text = $text
    $global:plainAnchor = $global:plainAnchor + 1
    return $newPara
}

function Insert-BoldParagraph([string]$text) {
    $target = $d.Paragraphs.Item($global:plainAnchor)
    $target.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($global:plainAnchor)
    $newPara.Range.Text = $text
    $newPara.Range.Font.Bold = -1
    $global:plainAnchor = $global:plainAnchor + 1
    return $newPara
}

# "Sun Sep 16" + " 12:47:22 PDT 2017" - authored as two separate runs
# in the source document, so force a run split at the boundary.
$dateLine = Insert-PlainParagraph("Sun Sep 16 12:47:22 PDT 2017")
$splitStart = $dateLine.Range.Start + 10
$splitEnd = $dateLine.Range.End - 1
$subRange = $d.Range($splitStart, $splitEnd)
$subRange.Font.Bold = -1
$subRange.Font.Bold = 0

Insert-PlainParagraph("Person Name" + "`t" + "`t" + "`t" + "`t" + "- M") | Out-Null
Insert-PlainParagraph("---------------------------------------------------------------") | Out-Null
Insert-PlainParagraph("Item Name" + "`t" + "`t" + "`t" + "`t" + "- CARROT") | Out-Null
Insert-PlainParagraph("Number of Pockets" + "`t" + "`t" + "`t" + "- 4") | Out-Null
Insert-PlainParagraph("Number of KGs" + "`t" + "`t" + "`t" + "- 372") | Out-Null
Insert-PlainParagraph("Rate" + "`t" + "`t" + "`t" + "`t" + "`t" + "- 20") | Out-Null
Insert-PlainParagraph("Transport & Miscellaneous" + "`t" + "- 40") | Out-Null
Insert-PlainParagraph("Total Price" + "`t" + "`t" + "`t" + "`t" + "- 7480.0") | Out-Null
Insert-BoldParagraph("Amount balance" + "`t" + "`t" + "`t" + "- 145817.0") | Out-Null
Insert-PlainParagraph("") | Out-Null
Insert-BoldParagraph("") | Out-Null

Write-Output "done"
